$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: W(kg) for Lophius budegassa corrected to 0
$ws.Range("G8").Value = 0

# Rows 31-46: updated Raising Factor (RF) value in column I
$ws.Range("I31:I46").Value = 17.20837209302325

# Row 42 and 45: Numb corrected to -1
$ws.Range("H42").Value = -1
$ws.Range("H45").Value = -1
